# Apply cryptos list update (Fri Dec  8 11:27:49 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold numeric-looking text (e.g. "43.548.38") that Excel would
# otherwise auto-convert to a number; force them to Text format first so the
# literal string is preserved exactly, one cell at a time (Range() with a comma-
# separated multi-area address only honors the first area for NumberFormat).
$textCells = @(
    "D2", "D3", "D5", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D36", "D37", "D38", "D40", "D43", "D44", "D45", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values in sheet order.
$ws.Range("D2").Value = "43.548.38"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.367.60"
$ws.Range("E3").Value = "  +5.58%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "234.05"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "72.63"
$ws.Range("E7").Value = "  +13.74%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").Value = "  +14.95%  "
$ws.Range("D10").Value = "0.0984"
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("D11").Value = "27.44"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "2.713.19"
$ws.Range("E12").Value = "  +5.35%  "
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("D14").Value = "16.40"
$ws.Range("E14").Value = "  +7.58%  "
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").Value = "0.871"
$ws.Range("E16").Value = "  +5.90%  "
$ws.Range("D17").Value = "2.362.17"
$ws.Range("E17").Value = "  +5.29%  "
$ws.Range("D18").Value = "43.447.13"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +4.48%  "
$ws.Range("D20").Value = "75.47"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("D21").Value = "6.37"
$ws.Range("E21").Value = "  +5.08%  "
$ws.Range("D22").Value = "251.56"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").Value = "3.81"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "2.47"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("D26").Value = "10.14"
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "22.69"
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("D29").Value = "172.81"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "1.54"
$ws.Range("E30").Value = "  +9.38%  "
$ws.Range("D31").Value = "0.134"
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").Value = "5.04"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("D36").Value = "3.77"
$ws.Range("E36").Value = "  +5.14%  "
$ws.Range("D37").Value = "6.66"
$ws.Range("E37").Value = "  +5.64%  "
$ws.Range("D38").Value = "2.45"
$ws.Range("E38").Value = "  +8.18%  "
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("D40").Value = "19.50"
$ws.Range("E40").Value = "  +14.58%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +3.90%  "
$ws.Range("D43").Value = "100.18"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "4.53"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  +10.76%  "
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").Value = "0.0960"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").Value = "1.446.54"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "2.589.93"
$ws.Range("E49").Value = "  +5.69%  "
$ws.Range("D50").Value = "2.77"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.173"
$ws.Range("E51").Value = "  +7.48%  "

Write-Output "Applied 93 cell updates"
